$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("feedstock_to_commodity")
$ws3 = $wb.Worksheets.Item("portfolio_metadata")

# -----------------------------------------------------------------------
# Fix biomethane: the buildout ramp for "Biomass -> Biomethane" (row 13)
# was wired to the Animal Manure feedstock (row 14) one row too early, and
# Animal Manure (row 14) / Diverted Organic Waste (row 15) were in turn
# pulling their 2045 quantities from the wrong Sheet1 rows. Biomethane
# itself (row 13) has no separate feedstock input of its own -- it is
# produced from the two feedstocks below it -- so its ramp now tracks the
# Animal Manure row's target (AA14) instead of its own, Animal Manure
# tracks Diverted Organic Waste's target (AA15), and Diverted Organic
# Waste tracks the next (empty) feedstock's target (AA16). AA13 becomes a
# plain static 0 instead of re-pulling Sheet1!C16.
# -----------------------------------------------------------------------

$ws2.Range("AA14").Formula = "=Sheet1!C16*1000"
$ws2.Range("AA15").Formula = "=Sheet1!C17*1000"

$ws2.Range("G13").Formula = "=AA13/21"
$ws2.Range("H13").Formula = "=`$AA14/21+G13"
$ws2.Range("I13").Formula = "=`$AA14/21+H13"
$ws2.Range("J13").Formula = "=`$AA14/21+I13"
$ws2.Range("K13").Formula = "=`$AA14/21+J13"
$ws2.Range("L13").Formula = "=`$AA14/21+K13"
$ws2.Range("M13").Formula = "=`$AA14/21+L13"
$ws2.Range("N13").Formula = "=`$AA14/21+M13"
$ws2.Range("O13").Formula = "=`$AA14/21+N13"
$ws2.Range("P13").Formula = "=`$AA14/21+O13"
$ws2.Range("Q13").Formula = "=`$AA14/21+P13"
$ws2.Range("R13").Formula = "=`$AA14/21+Q13"
$ws2.Range("S13").Formula = "=`$AA14/21+R13"
$ws2.Range("T13").Formula = "=`$AA14/21+S13"
$ws2.Range("U13").Formula = "=`$AA14/21+T13"
$ws2.Range("V13").Formula = "=`$AA14/21+U13"
$ws2.Range("W13").Formula = "=`$AA14/21+V13"
$ws2.Range("X13").Formula = "=`$AA14/21+W13"
$ws2.Range("Y13").Formula = "=`$AA14/21+X13"
$ws2.Range("Z13").Formula = "=`$AA14/21+Y13"

$ws2.Range("G14").Formula = "=AA14/21"
$ws2.Range("H14").Formula = "=`$AA15/21+G14"
$ws2.Range("I14").Formula = "=`$AA15/21+H14"
$ws2.Range("J14").Formula = "=`$AA15/21+I14"
$ws2.Range("K14").Formula = "=`$AA15/21+J14"
$ws2.Range("L14").Formula = "=`$AA15/21+K14"
$ws2.Range("M14").Formula = "=`$AA15/21+L14"
$ws2.Range("N14").Formula = "=`$AA15/21+M14"
$ws2.Range("O14").Formula = "=`$AA15/21+N14"
$ws2.Range("P14").Formula = "=`$AA15/21+O14"
$ws2.Range("Q14").Formula = "=`$AA15/21+P14"
$ws2.Range("R14").Formula = "=`$AA15/21+Q14"
$ws2.Range("S14").Formula = "=`$AA15/21+R14"
$ws2.Range("T14").Formula = "=`$AA15/21+S14"
$ws2.Range("U14").Formula = "=`$AA15/21+T14"
$ws2.Range("V14").Formula = "=`$AA15/21+U14"
$ws2.Range("W14").Formula = "=`$AA15/21+V14"
$ws2.Range("X14").Formula = "=`$AA15/21+W14"
$ws2.Range("Y14").Formula = "=`$AA15/21+X14"
$ws2.Range("Z14").Formula = "=`$AA15/21+Y14"

$ws2.Range("G15").Formula = "=AA15/21"
$ws2.Range("H15").Formula = "=`$AA16/21+G15"
$ws2.Range("I15").Formula = "=`$AA16/21+H15"
$ws2.Range("J15").Formula = "=`$AA16/21+I15"
$ws2.Range("K15").Formula = "=`$AA16/21+J15"
$ws2.Range("L15").Formula = "=`$AA16/21+K15"
$ws2.Range("M15").Formula = "=`$AA16/21+L15"
$ws2.Range("N15").Formula = "=`$AA16/21+M15"
$ws2.Range("O15").Formula = "=`$AA16/21+N15"
$ws2.Range("P15").Formula = "=`$AA16/21+O15"
$ws2.Range("Q15").Formula = "=`$AA16/21+P15"
$ws2.Range("R15").Formula = "=`$AA16/21+Q15"
$ws2.Range("S15").Formula = "=`$AA16/21+R15"
$ws2.Range("T15").Formula = "=`$AA16/21+S15"
$ws2.Range("U15").Formula = "=`$AA16/21+T15"
$ws2.Range("V15").Formula = "=`$AA16/21+U15"
$ws2.Range("W15").Formula = "=`$AA16/21+V15"
$ws2.Range("X15").Formula = "=`$AA16/21+W15"
$ws2.Range("Y15").Formula = "=`$AA16/21+X15"
$ws2.Range("Z15").Formula = "=`$AA16/21+Y15"

$ws2.Range("AA13").ClearFormats()
$ws2.Range("AA13").Value = 0

# -----------------------------------------------------------------------
# View-state: clear the (no-op) alignment formatting that had been left on
# portfolio_metadata!B5, then restore each sheet's last selection and make
# feedstock_to_commodity the active tab again.
# -----------------------------------------------------------------------

$ws3.Range("B5").ClearFormats()

$ws1.Activate()
$ws1.Range("B16").Select()

$ws3.Activate()
$ws3.Range("A1:B5").Select()

$ws2.Activate()
$ws2.Range("I31").Select()
